# Weekly refresh: a new "Apio" (celery) price pair (Primera/Segunda) for
# Terminal La Palmera de La Serena is inserted at row 435-436, pushing all
# subsequent rows (formerly 435-460) down by two rows to 437-462.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 435 (shifts 435:460 -> 437:462).
$ws.Rows("435:436").Insert()

# New row 435 - "Primera" quality entry for the newly added date.
$ws.Range("A435").Value = 8
$ws.Range("B435").Value = "Terminal La Palmera de La Serena"
$ws.Range("C435").Value = "Coquimbo"
$ws.Range("D435").Value = 44783
$ws.Range("E435").Value = 4
$ws.Range("F435").Value = 100112017
$ws.Range("G435").Value = "Apio"
$ws.Range("H435").Value = "Americana (o)"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 2000
$ws.Range("K435").Value = 7500
$ws.Range("L435").Value = 8000
$ws.Range("M435").Value = 7750
$ws.Range("N435").Value = "`$/docena de matas"
$ws.Range("O435").Value = "Provincia del Elquí"
$ws.Range("P435").Value = 1292
$ws.Range("Q435").Value = 6
$ws.Range("R435").Value = "Hortaliza"

# New row 436 - "Segunda" quality entry for the newly added date.
$ws.Range("A436").Value = 8
$ws.Range("B436").Value = "Terminal La Palmera de La Serena"
$ws.Range("C436").Value = "Coquimbo"
$ws.Range("D436").Value = 44783
$ws.Range("E436").Value = 4
$ws.Range("F436").Value = 100112017
$ws.Range("G436").Value = "Apio"
$ws.Range("H436").Value = "Americana (o)"
$ws.Range("I436").Value = "Segunda"
$ws.Range("J436").Value = 1200
$ws.Range("K436").Value = 6500
$ws.Range("L436").Value = 7000
$ws.Range("M436").Value = 6750
$ws.Range("N436").Value = "`$/docena de matas"
$ws.Range("O436").Value = "Provincia del Elquí"
$ws.Range("P436").Value = 1125
$ws.Range("Q436").Value = 6
$ws.Range("R436").Value = "Hortaliza"

# Ensure the date cells keep the workbook's date number format (style index 2
# in the original file), matching the other D-column cells.
$ws.Range("D435").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D436").NumberFormat = "YYYY-MM-DD HH:MM:SS"
